$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting the existing row 7 (and below) down to row 8
$ws.Rows("7:7").Insert()

# Fill the new row 7 with the new weekly record (same pattern as other rows, new date)
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C7").Value = "Los Lagos"
$ws.Range("D7").Value = 44505
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = "Otros"
$ws.Range("I7").Value = 100107002
$ws.Range("J7").Value = "Chirimoya"
$ws.Range("K7").Value = "Cultivar IV Región"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("Q7").Value = "$/bandeja 8 kilos"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 2438
$ws.Range("T7").Value = 8
